$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing values
$ws.Range("Q55").Value = 0
$ws.Range("Q59").Value = 0
$ws.Range("O344").Value = 2
$ws.Range("R346").Value = 0
$ws.Range("R347").Value = 0

# Append new weekly rows 348-362
# Row 348
$ws.Range("A348").Value = 45474
$ws.Range("A348").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B348").Value = 590.1500244140625
$ws.Range("C348").Value = 608.7999877929688
$ws.Range("D348").Value = 586.5
$ws.Range("E348").Value = 607.3499755859375
$ws.Range("F348").Value = 607.3499755859375
$ws.Range("G348").Value = 20476250
$ws.Range("H348").Value = 2024
$ws.Range("I348").Value = 7
$ws.Range("J348").Value = 1
$ws.Range("K348").Value = 0
$ws.Range("L348").Value = 0
$ws.Range("M348").Value = 0
$ws.Range("N348").Value = 27
$ws.Range("O348").Value = 0
$ws.Range("P348").Value = 0
$ws.Range("Q348").Value = 0

# Row 349
$ws.Range("A349").Value = 45481
$ws.Range("A349").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B349").Value = 609
$ws.Range("C349").Value = 640.75
$ws.Range("D349").Value = 608.5999755859375
$ws.Range("E349").Value = 635.5499877929688
$ws.Range("F349").Value = 635.5499877929688
$ws.Range("G349").Value = 24205618
$ws.Range("H349").Value = 2024
$ws.Range("I349").Value = 7
$ws.Range("J349").Value = 8
$ws.Range("K349").Value = 0
$ws.Range("L349").Value = 0
$ws.Range("M349").Value = 0
$ws.Range("N349").Value = 28
$ws.Range("O349").Value = 0
$ws.Range("P349").Value = 0
$ws.Range("Q349").Value = 0

# Row 350
$ws.Range("A350").Value = 45488
$ws.Range("A350").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B350").Value = 635.5499877929688
$ws.Range("C350").Value = 655
$ws.Range("D350").Value = 630.5
$ws.Range("E350").Value = 635.9000244140625
$ws.Range("F350").Value = 635.9000244140625
$ws.Range("G350").Value = 22590780
$ws.Range("H350").Value = 2024
$ws.Range("I350").Value = 7
$ws.Range("J350").Value = 15
$ws.Range("K350").Value = 0
$ws.Range("L350").Value = 0
$ws.Range("M350").Value = 0
$ws.Range("N350").Value = 29
$ws.Range("O350").Value = 0
$ws.Range("P350").Value = 0
$ws.Range("Q350").Value = 1

# Row 351
$ws.Range("A351").Value = 45495
$ws.Range("A351").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B351").Value = 631.9000244140625
$ws.Range("C351").Value = 709.7999877929688
$ws.Range("D351").Value = 623.25
$ws.Range("E351").Value = 703.3499755859375
$ws.Range("F351").Value = 703.3499755859375
$ws.Range("G351").Value = 43634651
$ws.Range("H351").Value = 2024
$ws.Range("I351").Value = 7
$ws.Range("J351").Value = 22
$ws.Range("K351").Value = 0
$ws.Range("L351").Value = 0
$ws.Range("M351").Value = 0
$ws.Range("N351").Value = 30
$ws.Range("O351").Value = 0
$ws.Range("P351").Value = 0
$ws.Range("Q351").Value = 0

# Row 352
$ws.Range("A352").Value = 45502
$ws.Range("A352").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B352").Value = 703.5499877929688
$ws.Range("C352").Value = 722.5
$ws.Range("D352").Value = 683.9000244140625
$ws.Range("E352").Value = 707.4000244140625
$ws.Range("F352").Value = 707.4000244140625
$ws.Range("G352").Value = 24279985
$ws.Range("H352").Value = 2024
$ws.Range("I352").Value = 7
$ws.Range("J352").Value = 29
$ws.Range("K352").Value = 0
$ws.Range("L352").Value = 0
$ws.Range("M352").Value = 0
$ws.Range("N352").Value = 31
$ws.Range("O352").Value = 0
$ws.Range("P352").Value = 0
$ws.Range("Q352").Value = 0

# Row 353
$ws.Range("A353").Value = 45509
$ws.Range("A353").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B353").Value = 693
$ws.Range("C353").Value = 719
$ws.Range("D353").Value = 677.2999877929688
$ws.Range("E353").Value = 702.4000244140625
$ws.Range("F353").Value = 702.4000244140625
$ws.Range("G353").Value = 29389482
$ws.Range("H353").Value = 2024
$ws.Range("I353").Value = 8
$ws.Range("J353").Value = 5
$ws.Range("K353").Value = 0
$ws.Range("L353").Value = 0
$ws.Range("M353").Value = 0
$ws.Range("N353").Value = 32
$ws.Range("O353").Value = 0
$ws.Range("P353").Value = 0
$ws.Range("Q353").Value = 0

# Row 354
$ws.Range("A354").Value = 45516
$ws.Range("A354").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B354").Value = 698.5999755859375
$ws.Range("C354").Value = 713.6500244140625
$ws.Range("D354").Value = 664.6500244140625
$ws.Range("E354").Value = 689.4000244140625
$ws.Range("F354").Value = 689.4000244140625
$ws.Range("G354").Value = 17011347
$ws.Range("H354").Value = 2024
$ws.Range("I354").Value = 8
$ws.Range("J354").Value = 12
$ws.Range("K354").Value = 0
$ws.Range("L354").Value = 0
$ws.Range("M354").Value = 0
$ws.Range("N354").Value = 33
$ws.Range("O354").Value = 0
$ws.Range("P354").Value = 0
$ws.Range("Q354").Value = 0

# Row 355
$ws.Range("A355").Value = 45523
$ws.Range("A355").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B355").Value = 689.4000244140625
$ws.Range("C355").Value = 730.5
$ws.Range("D355").Value = 678.5
$ws.Range("E355").Value = 729
$ws.Range("F355").Value = 729
$ws.Range("G355").Value = 23156540
$ws.Range("H355").Value = 2024
$ws.Range("I355").Value = 8
$ws.Range("J355").Value = 19
$ws.Range("K355").Value = 0
$ws.Range("L355").Value = 0
$ws.Range("M355").Value = 0
$ws.Range("N355").Value = 34
$ws.Range("O355").Value = 0
$ws.Range("P355").Value = 0
$ws.Range("Q355").Value = 0

# Row 356
$ws.Range("A356").Value = 45530
$ws.Range("A356").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B356").Value = 734
$ws.Range("C356").Value = 748.5
$ws.Range("D356").Value = 719.25
$ws.Range("E356").Value = 738.7000122070312
$ws.Range("F356").Value = 738.7000122070312
$ws.Range("G356").Value = 20945594
$ws.Range("H356").Value = 2024
$ws.Range("I356").Value = 8
$ws.Range("J356").Value = 26
$ws.Range("K356").Value = 0
$ws.Range("L356").Value = 0
$ws.Range("M356").Value = 0
$ws.Range("N356").Value = 35
$ws.Range("O356").Value = 0
$ws.Range("P356").Value = 0
$ws.Range("Q356").Value = 0

# Row 357
$ws.Range("A357").Value = 45537
$ws.Range("A357").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B357").Value = 743.1500244140625
$ws.Range("C357").Value = 761.2000122070312
$ws.Range("D357").Value = 736.0499877929688
$ws.Range("E357").Value = 738.4000244140625
$ws.Range("F357").Value = 738.4000244140625
$ws.Range("G357").Value = 16898676
$ws.Range("H357").Value = 2024
$ws.Range("I357").Value = 9
$ws.Range("J357").Value = 2
$ws.Range("K357").Value = 0
$ws.Range("L357").Value = 0
$ws.Range("M357").Value = 0
$ws.Range("N357").Value = 36
$ws.Range("O357").Value = 1
$ws.Range("P357").Value = 0
$ws.Range("Q357").Value = 0

# Row 358
$ws.Range("A358").Value = 45544
$ws.Range("A358").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B358").Value = 735.0999755859375
$ws.Range("C358").Value = 749
$ws.Range("D358").Value = 696
$ws.Range("E358").Value = 705.4000244140625
$ws.Range("F358").Value = 705.4000244140625
$ws.Range("G358").Value = 19980366
$ws.Range("H358").Value = 2024
$ws.Range("I358").Value = 9
$ws.Range("J358").Value = 9
$ws.Range("K358").Value = 0
$ws.Range("L358").Value = 0
$ws.Range("M358").Value = 0
$ws.Range("N358").Value = 37
$ws.Range("O358").Value = 0
$ws.Range("P358").Value = 0
$ws.Range("Q358").Value = 0

# Row 359
$ws.Range("A359").Value = 45551
$ws.Range("A359").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B359").Value = 708.8499755859375
$ws.Range("C359").Value = 721
$ws.Range("D359").Value = 691.5499877929688
$ws.Range("E359").Value = 711.0999755859375
$ws.Range("F359").Value = 711.0999755859375
$ws.Range("G359").Value = 17361938
$ws.Range("H359").Value = 2024
$ws.Range("I359").Value = 9
$ws.Range("J359").Value = 16
$ws.Range("K359").Value = 0
$ws.Range("L359").Value = 0
$ws.Range("M359").Value = 0
$ws.Range("N359").Value = 38
$ws.Range("O359").Value = 0
$ws.Range("P359").Value = 0
$ws.Range("Q359").Value = 0

# Row 360
$ws.Range("A360").Value = 45558
$ws.Range("A360").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B360").Value = 718.8499755859375
$ws.Range("C360").Value = 740.6500244140625
$ws.Range("D360").Value = 710.0999755859375
$ws.Range("E360").Value = 730.2000122070312
$ws.Range("F360").Value = 730.2000122070312
$ws.Range("G360").Value = 14617652
$ws.Range("H360").Value = 2024
$ws.Range("I360").Value = 9
$ws.Range("J360").Value = 23
$ws.Range("K360").Value = 0
$ws.Range("L360").Value = 0
$ws.Range("M360").Value = 0
$ws.Range("N360").Value = 39
$ws.Range("O360").Value = 0
$ws.Range("P360").Value = 0
$ws.Range("Q360").Value = 0

# Row 361
$ws.Range("A361").Value = 45565
$ws.Range("A361").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B361").Value = 735
$ws.Range("C361").Value = 736.3499755859375
$ws.Range("D361").Value = 696.5
$ws.Range("E361").Value = 708.7999877929688
$ws.Range("F361").Value = 708.7999877929688
$ws.Range("G361").Value = 9635468
$ws.Range("H361").Value = 2024
$ws.Range("I361").Value = 9
$ws.Range("J361").Value = 30
$ws.Range("K361").Value = 0
$ws.Range("L361").Value = 0
$ws.Range("M361").Value = 0
$ws.Range("N361").Value = 40
$ws.Range("O361").Value = 0
$ws.Range("P361").Value = 0
$ws.Range("Q361").Value = 0

# Row 362
$ws.Range("A362").Value = 45572
$ws.Range("A362").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B362").Value = 711.6500244140625
$ws.Range("C362").Value = 726.2000122070312
$ws.Range("D362").Value = 699.7999877929688
$ws.Range("E362").Value = 723.8499755859375
$ws.Range("F362").Value = 723.8499755859375
$ws.Range("G362").Value = 13218801
$ws.Range("H362").Value = 2024
$ws.Range("I362").Value = 10
$ws.Range("J362").Value = 7
$ws.Range("K362").Value = 0
$ws.Range("L362").Value = 0
$ws.Range("M362").Value = 0
$ws.Range("N362").Value = 41
$ws.Range("O362").Value = 0
$ws.Range("P362").Value = 0
$ws.Range("Q362").Value = 0
$ws.Range("R362").Value = 0
